# Update the "Förändrad" (Changed) date column (C) for data rows 2-12
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
